# BaseRules.xlsx — update the truth-table "FS" output column (column B)
# for several rows, update E13, and move the active selection to D25.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B updates (rows 6-9, 14-17, 22-25): each block's values bump by 1
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 3

$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 3

$ws.Range("B22").Value = 2
$ws.Range("B23").Value = 2
$ws.Range("B24").Value = 3
$ws.Range("B25").Value = 3

# Column E update (row 13)
$ws.Range("E13").Value = 1

# Move the selection / scroll position to match the saved view state
[void]$ws.Range("D25").Select()
